$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "Controlar la consistencia y la integridad de los datos: "
# The trailing space after the (underlined) label is split off into its
# own plain run so it no longer carries the underline formatting.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Controlar la consistencia y la integridad de los datos: ")
if ($found) {
    $spaceStart = $rng.End - 1
    $spaceEnd = $rng.End
    $space = $d.Range($spaceStart, $spaceEnd)
    $space.Font.Underline = 0
}

# ---------------------------------------------------------------------
# Change 2: " de base de datos: Ofrecer copias de seguridad, ..."
# The trailing space after the (underlined) label is split off into its
# own plain run, and the leading "O" of "Ofrecer" is lower-cased to "o"
# and folded into the underlined label run, leaving "frecer..." as a
# separate plain run.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(" de base de datos: Ofrecer")
if ($found) {
    $matchEnd = $rng.End
    $oStart = $matchEnd - 7      # length of "Ofrecer"
    $oEnd = $oStart + 1
    $spaceStart = $oStart - 1
    $spaceEnd = $oStart

    $space = $d.Range($spaceStart, $spaceEnd)
    $space.Font.Underline = 0

    $oRange = $d.Range($oStart, $oEnd)
    $oRange.Text = "o"
    $oRange2 = $d.Range($oStart, $oStart + 1)
    $oRange2.Font.Underline = 1
}

# ---------------------------------------------------------------------
# Change 3: "Controlar la concurrencia: "
# Same trailing-space split as change 1.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Controlar la concurrencia: ")
if ($found) {
    $spaceStart = $rng.End - 1
    $spaceEnd = $rng.End
    $space = $d.Range($spaceStart, $spaceEnd)
    $space.Font.Underline = 0
}
